$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 9770.808000000001
$ws.Range("I33").Value = 11146.762
$ws.Range("J33").Value = 3991.8
$ws.Range("K33").Value = 11146.762
$ws.Range("L33").Value = 3991.8
$ws.Range("M33").Value = -10917.762
$ws.Range("N33").Value = -4449.8
$ws.Range("H88").Value = 2628.5833
$ws.Range("I88").Value = 1690
$ws.Range("J88").Value = 3097.875
$ws.Range("K88").Value = 1690
$ws.Range("L88").Value = 3097.875
$ws.Range("M88").Value = -1284
$ws.Range("N88").Value = -3909.875
$ws.Range("H91").Value = 2628.5833
$ws.Range("I91").Value = 1690
$ws.Range("J91").Value = 3097.875
$ws.Range("K91").Value = 1690
$ws.Range("L91").Value = 3097.875
$ws.Range("M91").Value = -286
$ws.Range("N91").Value = -5905.875
$ws.Range("H125").Value = 3953.25
$ws.Range("I125").Value = 3387.3333
$ws.Range("K125").Value = 30485.9997
$ws.Range("M125").Value = -28025.9997
$ws.Range("H135").Value = 3271.5881
$ws.Range("I135").Value = 3432.375
$ws.Range("J135").Value = 699
$ws.Range("K135").Value = 30891.375
$ws.Range("L135").Value = 6291
$ws.Range("M135").Value = -28356.375
$ws.Range("N135").Value = -11361
$ws.Range("H137").Value = 6047.4053
$ws.Range("I137").Value = 3887.0356
$ws.Range("K137").Value = 11661.1068
$ws.Range("M137").Value = -9111.106800000001
$ws.Range("H138").Value = 34485900
$ws.Range("I138").Value = 1555.0769
$ws.Range("J138").Value = 62504430
$ws.Range("K138").Value = 4665.2307
$ws.Range("L138").Value = 187513290
$ws.Range("M138").Value = 474.7692999999999
$ws.Range("N138").Value = -187523570

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1676.5714
$ws.Range("I2").Value = 1148.2
$ws.Range("K2").Value = 1148.2
$ws.Range("M2").Value = -1035.2
$ws.Range("H32").Value = 2778.0417
$ws.Range("I32").Value = 2767.3838
$ws.Range("K32").Value = 2767.3838
$ws.Range("M32").Value = -2480.3838
$ws.Range("H63").Value = 2425
$ws.Range("I63").Value = 2425
$ws.Range("K63").Value = 2425
$ws.Range("M63").Value = -1739
$ws.Range("H66").Value = 2425
$ws.Range("I66").Value = 2425
$ws.Range("K66").Value = 12125
$ws.Range("M66").Value = -8693
$ws.Range("H74").Value = 38008.547
$ws.Range("I74").Value = 40436.723
$ws.Range("K74").Value = 40436.723
$ws.Range("M74").Value = -39562.723
$ws.Range("H77").Value = 38008.547
$ws.Range("I77").Value = 40436.723
$ws.Range("K77").Value = 202183.615
$ws.Range("M77").Value = -197815.615
$ws.Range("H116").Value = 1676.5714
$ws.Range("I116").Value = 1148.2
$ws.Range("K116").Value = 1148.2
$ws.Range("M116").Value = 1145.8
$ws.Range("H132").Value = 88067.42999999999
$ws.Range("I132").Value = 2071.7058
$ws.Range("J132").Value = 220969.9
$ws.Range("K132").Value = 6215.117400000001
$ws.Range("L132").Value = 662909.7
$ws.Range("M132").Value = -3685.117400000001
$ws.Range("N132").Value = -667969.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1676.5714
$ws.Range("I3").Value = 1148.2
$ws.Range("K3").Value = 1148.2
$ws.Range("M3").Value = -1034.2
$ws.Range("H20").Value = 2047.2222
$ws.Range("I20").Value = 2788.8
$ws.Range("J20").Value = 1120.25
$ws.Range("K20").Value = 2788.8
$ws.Range("L20").Value = 1120.25
$ws.Range("M20").Value = -2541.8
$ws.Range("N20").Value = -1614.25
$ws.Range("H69").Value = 120000
$ws.Range("I69").Value = 100000
$ws.Range("K69").Value = 100000
$ws.Range("M69").Value = -99189
$ws.Range("H72").Value = 120000
$ws.Range("I72").Value = 100000
$ws.Range("K72").Value = 300000
$ws.Range("M72").Value = -295944
$ws.Range("H86").Value = 18483.428
$ws.Range("I86").Value = 11114.429
$ws.Range("J86").Value = 33221.43
$ws.Range("K86").Value = 11114.429
$ws.Range("L86").Value = 33221.43
$ws.Range("M86").Value = -9991.429
$ws.Range("N86").Value = -35467.43
$ws.Range("H89").Value = 18483.428
$ws.Range("I89").Value = 11114.429
$ws.Range("J89").Value = 33221.43
$ws.Range("K89").Value = 55572.145
$ws.Range("L89").Value = 166107.15
$ws.Range("M89").Value = -49956.145
$ws.Range("N89").Value = -177339.15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 252.30435
$ws.Range("I7").Value = 121
$ws.Range("J7").Value = 353.30768
$ws.Range("K7").Value = 121
$ws.Range("L7").Value = 353.30768
$ws.Range("M7").Value = -8
$ws.Range("N7").Value = -579.30768

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H122").Value = 1350.0416
$ws.Range("I122").Value = 1499
$ws.Range("J122").Value = 1343.5652
$ws.Range("K122").Value = 13491
$ws.Range("L122").Value = 12092.0868
$ws.Range("M122").Value = -11041
$ws.Range("N122").Value = -16992.0868
$ws.Range("H131").Value = 1363
$ws.Range("I131").Value = 943
$ws.Range("J131").Value = 1514.2
$ws.Range("K131").Value = 2829
$ws.Range("L131").Value = 4542.6
$ws.Range("M131").Value = 2211
$ws.Range("N131").Value = -14622.6
$ws.Range("H134").Value = 2955.3125
$ws.Range("I134").Value = 1234.6428
$ws.Range("K134").Value = 3703.9284
$ws.Range("M134").Value = 1366.0716
$ws.Range("H140").Value = 795.8889
$ws.Range("I140").Value = 795.8889
$ws.Range("K140").Value = 2387.6667
$ws.Range("M140").Value = 2792.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6007.8887
$ws.Range("I70").Value = 4999.077
$ws.Range("J70").Value = 8630.799999999999
$ws.Range("K70").Value = 4999.077
$ws.Range("L70").Value = 8630.799999999999
$ws.Range("M70").Value = -4729.077
$ws.Range("N70").Value = -9170.799999999999
$ws.Range("H73").Value = 6007.8887
$ws.Range("I73").Value = 4999.077
$ws.Range("J73").Value = 8630.799999999999
$ws.Range("K73").Value = 4999.077
$ws.Range("L73").Value = 8630.799999999999
$ws.Range("M73").Value = -4063.077
$ws.Range("N73").Value = -10502.8
$ws.Range("H102").Value = 8887.143
$ws.Range("I102").Value = 2868.3333
$ws.Range("K102").Value = 2868.3333
$ws.Range("M102").Value = -1246.3333
$ws.Range("H132").Value = 1191.619
$ws.Range("I132").Value = 1197.2
$ws.Range("K132").Value = 3591.6
$ws.Range("M132").Value = -1061.6
$ws.Range("H134").Value = 92326
$ws.Range("J134").Value = 92326
$ws.Range("L134").Value = 276978
$ws.Range("N134").Value = -282048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1884.2632
$ws.Range("I22").Value = 1755.7778
$ws.Range("J22").Value = 1999.9
$ws.Range("K22").Value = 1755.7778
$ws.Range("L22").Value = 1999.9
$ws.Range("M22").Value = -1460.7778
$ws.Range("N22").Value = -2589.9
$ws.Range("H27").Value = 1884.2632
$ws.Range("I27").Value = 1755.7778
$ws.Range("J27").Value = 1999.9
$ws.Range("K27").Value = 1755.7778
$ws.Range("L27").Value = 1999.9
$ws.Range("M27").Value = -1648.7778
$ws.Range("N27").Value = -2213.9
$ws.Range("H40").Value = 3992.5
$ws.Range("I40").Value = 3992.5
$ws.Range("K40").Value = 3992.5
$ws.Range("M40").Value = -3856.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 114779.445
$ws.Range("I122").Value = 146459.28
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 439377.84
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -436927.84
$ws.Range("N122").Value = -16600
$ws.Range("H136").Value = 6626.471
$ws.Range("I136").Value = 5427.9443
$ws.Range("J136").Value = 10671.5
$ws.Range("K136").Value = 16283.8329
$ws.Range("L136").Value = 32014.5
$ws.Range("M136").Value = -13733.8329
$ws.Range("N136").Value = -37114.5
